# Weekly price update: a new report row is inserted as row 20 (pushing the
# existing rows 20-99 down to 21-100, which is exactly what the diff shows -
# every row's data equals what used to sit one row above it, with a brand
# new record landing in the vacated row 20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 20; Excel shifts rows 20:99 down to 21:100
# and the sheet's used range / dimension grows from R99 to R100 automatically.
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new weekly record.
$ws.Range("A20").Value = 5
$ws.Range("B20").Value = "Macroferia Regional de Talca"
$ws.Range("C20").Value = "Maule"
$ws.Range("D20").Value = 44602
$ws.Range("E20").Value = 7
$ws.Range("F20").Value = 100112001
$ws.Range("G20").Value = "Berenjena"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 200
$ws.Range("K20").Value = 7000
$ws.Range("L20").Value = 7000
$ws.Range("M20").Value = 7000
$ws.Range("N20").Value = "`$/caja 50 unidades"
$ws.Range("O20").Value = "Región del Maule"
$ws.Range("P20").Value = 140
$ws.Range("Q20").Value = 50
$ws.Range("R20").Value = "Hortaliza"
